# Insert a new weekly record for "Betarraga" (Vega Monumental Concepción)
# as row 358, shifting all subsequent rows (old 358..459) down by one
# (they become 359..460). The workbook's used range grows from A1:R459
# to A1:R460.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 358 down one row.
$ws.Rows(358).Insert()

# Populate the newly-inserted row 358 with the new weekly record.
$ws.Range("A358").Value = 11
$ws.Range("B358").Value = "Vega Monumental Concepción"
$ws.Range("C358").Value = "Bíobío"
$ws.Range("D358").Value = 44985
$ws.Range("E358").Value = 8
$ws.Range("F358").Value = 100114014
$ws.Range("G358").Value = "Betarraga"
$ws.Range("H358").Value = "Sin especificar"
$ws.Range("I358").Value = "Primera"
$ws.Range("J358").Value = 650
$ws.Range("K358").Value = 600
$ws.Range("L358").Value = 650
$ws.Range("M358").Value = 627
$ws.Range("N358").Value = "$/paquete 5 unidades"
$ws.Range("O358").Value = "Región Metropolitana"
$ws.Range("P358").Value = 125
$ws.Range("Q358").Value = 5
$ws.Range("R358").Value = "Hortaliza"
